# edit.ps1 — applies the "Updated the definition document" change:
#  1. Splits the sentence "I've found this resource quite useful when working
#     with this algorithm." into four separate runs that together read
#     "I've found these resources quite useful when working with this algorithm."
#  2. Replaces the blank paragraph that follows the Stanford hyperlink
#     paragraph with a new hyperlink paragraph (linking to the
#     policyalmanac.org A* tutorial) followed by a brand-new empty paragraph.

$d = $word.ActiveDocument

function Get-ParagraphIndexContainingPosition($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pr = $d.Paragraphs.Item($i).Range
        if ($pos -ge $pr.Start -and $pos -lt $pr.End) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Part 1: split the "I've found this resource ..." run into four runs.
# ---------------------------------------------------------------------------

$targetText = "I've found this resource quite useful when working with this algorithm."
$searchRange = $d.Range(0, $d.Content.End)
$found = $searchRange.Find.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the sentence about the resource."
}

$paraIndex = Get-ParagraphIndexContainingPosition($searchRange.Start)
if ($paraIndex -eq -1) {
    throw "Could not locate the paragraph containing the resource sentence."
}

$paraRange = $d.Paragraphs.Item($paraIndex).Range
$xml = $paraRange.WordOpenXML

$runPattern = '(?s)(<w:r\b[^>]*>)(<w:rPr>.*?</w:rPr>)(<w:t[^>]*>)I''ve found this resource quite useful when working with this algorithm\.(</w:t>)(</w:r>)'
$runMatch = [System.Text.RegularExpressions.Regex]::Match($xml, $runPattern)
if (-not $runMatch.Success) {
    throw "Could not locate the run to split inside the paragraph's XML."
}

$openRunTag = $runMatch.Groups[1].Value
$runProps = $runMatch.Groups[2].Value

$run1 = $openRunTag + $runProps + "<w:t>I've found these</w:t></w:r>"
$run2 = "<w:r>" + $runProps + '<w:t xml:space="preserve"> resource</w:t></w:r>'
$run3 = "<w:r>" + $runProps + "<w:t>s</w:t></w:r>"
$run4 = "<w:r>" + $runProps + '<w:t xml:space="preserve"> quite useful when working with this algorithm.</w:t></w:r>'
$splitReplacement = $run1 + $run2 + $run3 + $run4

$xml = $xml.Substring(0, $runMatch.Index) + $splitReplacement + $xml.Substring($runMatch.Index + $runMatch.Length)

$paraRange.InsertXML($xml)

# ---------------------------------------------------------------------------
# Part 2: turn the blank paragraph after the Stanford hyperlink into a new
# hyperlink paragraph, followed by a fresh blank paragraph.
# ---------------------------------------------------------------------------

$stanfordParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*theory.stanford.edu*") {
        $stanfordParaIndex = $i
        break
    }
}
if ($stanfordParaIndex -eq -1) {
    throw "Could not locate the existing Stanford hyperlink paragraph."
}

$blankParaIndex = $stanfordParaIndex + 1
$blankParaRange = $d.Paragraphs.Item($blankParaIndex).Range
if ($blankParaRange.Text.Trim().Length -ne 0) {
    throw "Expected the paragraph following the Stanford hyperlink to be blank."
}

$xml2 = $blankParaRange.WordOpenXML
$newUrl = "http://www.policyalmanac.org/games/aStarTutorial.htm"

# Add a fresh relationship for the new external hyperlink.
$relsPattern = '(?s)(<pkg:part pkg:name="/word/_rels/document.xml.rels"[^>]*><pkg:xmlData><Relationships[^>]*>)(.*?)(</Relationships></pkg:xmlData></pkg:part>)'
$relsMatch = [System.Text.RegularExpressions.Regex]::Match($xml2, $relsPattern)
if (-not $relsMatch.Success) {
    throw "Could not locate the document relationships part."
}
$newRelationship = '<Relationship Id="rIdNewTutorialLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="' + $newUrl + '" TargetMode="External"/>'
$relsInner = $relsMatch.Groups[2].Value + $newRelationship
$xml2 = $xml2.Substring(0, $relsMatch.Groups[2].Index) + $relsInner + $xml2.Substring($relsMatch.Groups[2].Index + $relsMatch.Groups[2].Length)

# Replace the single blank <w:p> with the hyperlink paragraph + a new blank paragraph.
$bodyPattern = '(?s)<w:body>(<w:p\b[^>]*>)(<w:pPr>.*?</w:pPr>)?</w:p>'
$bodyMatch = [System.Text.RegularExpressions.Regex]::Match($xml2, $bodyPattern)
if (-not $bodyMatch.Success) {
    throw "Could not locate the blank paragraph's body markup."
}
$openParaTag = $bodyMatch.Groups[1].Value
$newBody = '<w:body>' + $openParaTag + '<w:hyperlink r:id="rIdNewTutorialLink" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>' + $newUrl + '</w:t></w:r></w:hyperlink></w:p><w:p/>'
$xml2 = $xml2.Substring(0, $bodyMatch.Index) + $newBody + $xml2.Substring($bodyMatch.Index + $bodyMatch.Length)

$blankParaRange.InsertXML($xml2)

Write-Host "Done."
